# Casos_de_Testes_Grupo_9: add the 19/04/2018 test result row (D20 = "ok")
# and move the active selection down to it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New test case result appended at D20.
$ws.Range("D20").Value = "ok"

# Move the selection / view down to the newly added row, matching the
# author's final cursor position.
[void]$ws.Range("D20").Select()
$excel.ActiveWindow.ScrollRow = 11
$excel.ActiveWindow.ScrollColumn = 3
